$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("B41").Value = 'Stellar'
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D2").Value = '44.519.97'
$ws.Range("D3").Value = '2.428.13'
Set-TextValue "D4" '1.00'
Set-TextValue "D5" '313.99'
Set-TextValue "D6" '101.61'
Set-TextValue "D9" '0.512'
Set-TextValue "D10" '35.23'
Set-TextValue "D13" '18.75'
Set-TextValue "D14" '6.95'
$ws.Range("D15").Value = '2.809.86'
$ws.Range("D16").Value = '2.398.61'
Set-TextValue "D17" '0.837'
$ws.Range("D18").Value = '44.440.13'
Set-TextValue "D19" '12.44'
Set-TextValue "D20" '6.39'
$ws.Range("D21").Value = '0.0₃0907'
Set-TextValue "D22" '68.95'
Set-TextValue "D23" '241.09'
Set-TextValue "D25" '2.48'
Set-TextValue "D27" '25.16'
Set-TextValue "D30" '33.28'
Set-TextValue "D31" '48.44'
Set-TextValue "D32" '0.122'
Set-TextValue "D36" '0.0765'
Set-TextValue "D38" '4.52'
Set-TextValue "D40" '127.21'
Set-TextValue "D41" '0.109'
Set-TextValue "D42" '22.02'
$ws.Range("D45").Value = '1.946.93'
Set-TextValue "D47" '2.94'
Set-TextValue "D48" '9.74'
Set-TextValue "D50" '53.45'
Set-TextValue "D51" '73.76'
$ws.Range("E2").Value = '  +3.75%  '
$ws.Range("E3").Value = '  +2.75%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +3.99%  '
$ws.Range("E6").Value = '  +6.45%  '
$ws.Range("E7").Value = '  +1.60%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +5.37%  '
$ws.Range("E10").Value = '  +3.90%  '
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("E12").Value = '  +1.24%  '
$ws.Range("E13").Value = '  +2.77%  '
$ws.Range("E14").Value = '  +3.41%  '
$ws.Range("E15").Value = '  +2.88%  '
$ws.Range("E16").Value = '  +2.23%  '
$ws.Range("E17").Value = '  +5.17%  '
$ws.Range("E18").Value = '  +3.69%  '
$ws.Range("E19").Value = '  +3.27%  '
$ws.Range("E20").Value = '  +2.09%  '
$ws.Range("E21").Value = '  +2.46%  '
$ws.Range("E22").Value = '  +1.63%  '
$ws.Range("E23").Value = '  +2.59%  '
$ws.Range("E24").Value = '  +4.13%  '
$ws.Range("E25").Value = '  +2.48%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  +1.97%  '
$ws.Range("E28").Value = '  -4.23%  '
$ws.Range("E30").Value = '  +5.59%  '
$ws.Range("E31").Value = '  +1.15%  '
$ws.Range("E32").Value = '  +16.33%  '
$ws.Range("E33").Value = '  +12.63%  '
$ws.Range("E34").Value = '  +3.43%  '
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("E36").Value = '  +5.70%  '
$ws.Range("E37").Value = '  +2.75%  '
$ws.Range("E38").Value = '  +3.72%  '
$ws.Range("E39").Value = '  +3.98%  '
$ws.Range("E40").Value = '  +5.40%  '
$ws.Range("E41").Value = '  +0.93%  '
$ws.Range("E42").Value = '  +2.26%  '
$ws.Range("E43").Value = '  -5.58%  '
$ws.Range("E44").Value = '  +3.38%  '
$ws.Range("E45").Value = '  +0.70%  '
$ws.Range("E46").Value = '  +2.20%  '
$ws.Range("E47").Value = '  +8.48%  '
$ws.Range("E48").Value = '  +5.99%  '
$ws.Range("E49").Value = '  +11.48%  '
$ws.Range("E50").Value = '  +3.49%  '
$ws.Range("E51").Value = '  +2.42%  '
